$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("buffer")

# The "name" column (A) is being replaced by "id" -- i.e. column A is
# removed entirely and every other column shifts one place to the left
# (type, location, item, onhand).
$ws.Range("A1:D2").Value2 = $ws.Range("B1:E2").Value2
$ws.Columns.Item(5).Delete()

$null = $ws.Activate()
$null = $ws.Columns.Item(1).Select()
